$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 65490.97652451508
$ws.Range("B3").Value = 62526.16260984209
$ws.Range("B4").Value = 60367.94398585002
$ws.Range("B5").Value = 58598.65508066331
$ws.Range("B6").Value = 58928.72293532812
$ws.Range("B7").Value = 62707.00805384065
$ws.Range("B8").Value = 63780.04017997128
$ws.Range("B9").Value = 70341.91836723284
$ws.Range("B10").Value = 85318.32505643652
$ws.Range("B11").Value = 92332.91992569293
$ws.Range("B12").Value = 97055.31671802221
$ws.Range("B13").Value = 99389.32805701185
$ws.Range("B14").Value = 97464.16776457951
$ws.Range("B15").Value = 101461.6258321211
$ws.Range("B16").Value = 103092.1323862914
$ws.Range("B17").Value = 102442.8344496576
$ws.Range("B18").Value = 98107.23082194131
$ws.Range("B19").Value = 90031.66561785711
$ws.Range("B20").Value = 90333.16187821308
$ws.Range("B21").Value = 87548.79118056811
$ws.Range("B22").Value = 86118.20735377351
$ws.Range("B23").Value = 85962.51943495168
$ws.Range("B24").Value = 82330.94843868678
$ws.Range("B25").Value = 77291.86507444084
